# Update the "tree" benchmark worksheet: insert a new storage-cost table
# ("对应存储开销") between the existing build-time and query-time tables.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Centre-aligned style used for (almost) every cell on this sheet.
$xlCenter = -4108

function Set-Cell($row, $col, $value, $vAlign = $false) {
    $c = $ws.Cells.Item($row, $col)
    if ($null -ne $value) {
        $c.Value = $value
    }
    $c.HorizontalAlignment = $xlCenter
    if ($vAlign) {
        $c.VerticalAlignment = $xlCenter
    }
}

# Start from a clean sheet so the old row layout (rows 7-20) doesn't linger.
$ws.Cells.Clear()

# --- Row 1: section title ----------------------------------------------
Set-Cell 1 1 "索引建立（10次构建取平均）"

# --- Rows 2-5: build-time table ------------------------------------------
Set-Cell 2 1 "维度\数据量"
Set-Cell 2 2 5000
Set-Cell 2 3 10000
Set-Cell 2 4 15000
Set-Cell 2 5 20000
Set-Cell 2 6 25000
Set-Cell 2 7 $null

Set-Cell 3 1 "2 (4叉树)" $true
Set-Cell 3 2 190
Set-Cell 3 3 389
Set-Cell 3 4 568
Set-Cell 3 5 811
Set-Cell 3 6 1027
Set-Cell 3 7 "单位 ms"

Set-Cell 4 1 "3 (8叉树)" $true
Set-Cell 4 2 218
Set-Cell 4 3 450
Set-Cell 4 4 654
Set-Cell 4 5 885
Set-Cell 4 6 1092
Set-Cell 4 7 $null

Set-Cell 5 1 "4 (16叉树)" $true
Set-Cell 5 2 254
Set-Cell 5 3 506
Set-Cell 5 4 825
Set-Cell 5 5 1015
Set-Cell 5 6 1350
Set-Cell 5 7 $null

# --- Row 6: blank spacer row (keeps its centred formatting) --------------
Set-Cell 6 2 $null
Set-Cell 6 3 $null
Set-Cell 6 4 $null
Set-Cell 6 5 $null
Set-Cell 6 6 $null
Set-Cell 6 7 $null

# --- Row 8: new section title --------------------------------------------
Set-Cell 8 1 "对应存储开销" $true

# --- Rows 9-12: new storage-cost table ------------------------------------
Set-Cell 9 1 "维度\数据量"
Set-Cell 9 2 5000
Set-Cell 9 3 10000
Set-Cell 9 4 15000
Set-Cell 9 5 20000
Set-Cell 9 6 25000
Set-Cell 9 7 $null

Set-Cell 10 1 "2 (4叉树)" $true
Set-Cell 10 2 200000
Set-Cell 10 3 400000
Set-Cell 10 4 600000
Set-Cell 10 5 800000
Set-Cell 10 6 1000000
Set-Cell 10 7 "单位 bytes"

Set-Cell 11 1 "3 (8叉树)" $true
Set-Cell 11 2 200000
Set-Cell 11 3 400000
Set-Cell 11 4 600000
Set-Cell 11 5 800000
Set-Cell 11 6 1000000
Set-Cell 11 7 $null

Set-Cell 12 1 "4 (16叉树)" $true
Set-Cell 12 2 200000
Set-Cell 12 3 400000
Set-Cell 12 4 600000
Set-Cell 12 5 800000
Set-Cell 12 6 1000000
Set-Cell 12 7 $null

# --- Row 15: query-time section title (moved down from row 7) ------------
Set-Cell 15 1 "查询（随机查找100个已知节点的平均时间）" $true
Set-Cell 15 2 $null
Set-Cell 15 3 $null
Set-Cell 15 4 $null
Set-Cell 15 5 $null
Set-Cell 15 6 $null
Set-Cell 15 7 $null

# --- Rows 16-19: query-time table (moved down from rows 8-11) ------------
Set-Cell 16 1 "维度\数据量"
Set-Cell 16 2 5000
Set-Cell 16 3 10000
Set-Cell 16 4 15000
Set-Cell 16 5 20000
Set-Cell 16 6 25000
Set-Cell 16 7 $null

Set-Cell 17 1 "2 (4叉树)" $true
Set-Cell 17 2 566
Set-Cell 17 3 579
Set-Cell 17 4 551
Set-Cell 17 5 558
Set-Cell 17 6 539
Set-Cell 17 7 "单位 ms"

Set-Cell 18 1 "3 (8叉树)" $true
Set-Cell 18 2 544
Set-Cell 18 3 539
Set-Cell 18 4 561
Set-Cell 18 5 546
Set-Cell 18 6 551
Set-Cell 18 7 $null

Set-Cell 19 1 "4 (16叉树)" $true
Set-Cell 19 2 544
Set-Cell 19 3 546
Set-Cell 19 4 542
Set-Cell 19 5 549
Set-Cell 19 6 573
Set-Cell 19 7 $null

# --- New column 7 width (added alongside the new "单位 bytes" column) ----
$ws.Columns.Item(7).ColumnWidth = 14.25

# --- Selection, matching the commit's final cursor position --------------
[void]$ws.Range("K12").Select()
